$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3533.8928
$ws.Range("I32").Value = 2525.423
$ws.Range("K32").Value = 2525.423
$ws.Range("M32").Value = -2238.423

$ws.Range("H74").Value = 2696.9285
$ws.Range("I74").Value = 1917.2222
$ws.Range("K74").Value = 1917.2222
$ws.Range("M74").Value = -1043.2222

$ws.Range("H77").Value = 2696.9285
$ws.Range("I77").Value = 1917.2222
$ws.Range("K77").Value = 9586.110999999999
$ws.Range("M77").Value = -5218.110999999999

$ws.Range("H110").Value = 1457.4166
$ws.Range("I110").Value = 1457.4166
$ws.Range("K110").Value = 1457.4166
$ws.Range("M110").Value = 587.5834

$ws.Range("H132").Value = 1703.9656
$ws.Range("I132").Value = 1657.6786
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4973.0358
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2443.0358
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7016.6113
$ws.Range("I86").Value = 6480.875
$ws.Range("K86").Value = 6480.875
$ws.Range("M86").Value = -5357.875

$ws.Range("H89").Value = 7016.6113
$ws.Range("I89").Value = 6480.875
$ws.Range("K89").Value = 32404.375
$ws.Range("M89").Value = -26788.375

$ws.Range("H94").Value = 2051.5
$ws.Range("I94").Value = 1700.25
$ws.Range("K94").Value = 1700.25
$ws.Range("M94").Value = -1249.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 683.2222
$ws.Range("I86").Value = 649.8570999999999
$ws.Range("K86").Value = 1949.5713
$ws.Range("M86").Value = -763.5712999999998

$ws.Range("H89").Value = 683.2222
$ws.Range("I89").Value = 649.8570999999999
$ws.Range("K89").Value = 5848.7139
$ws.Range("M89").Value = 79.28610000000026

$ws.Range("H122").Value = 895.7059
$ws.Range("J122").Value = 1293.8
$ws.Range("L122").Value = 11644.2
$ws.Range("N122").Value = -16544.2

$ws.Range("H136").Value = 6035.154
$ws.Range("I136").Value = 3682.4546
$ws.Range("K136").Value = 11047.3638
$ws.Range("M136").Value = -5947.363799999999

$ws.Range("H139").Value = 10379.286
$ws.Range("I139").Value = 11862
$ws.Range("J139").Value = 9555.556
$ws.Range("K139").Value = 35586
$ws.Range("L139").Value = 28666.668
$ws.Range("M139").Value = -30446
$ws.Range("N139").Value = -38946.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 36474.5
$ws.Range("I26").Value = 39000
$ws.Range("J26").Value = 35632.668
$ws.Range("K26").Value = 39000
$ws.Range("L26").Value = 35632.668
$ws.Range("M26").Value = -38720
$ws.Range("N26").Value = -36192.668

$ws.Range("H50").Value = 36474.5
$ws.Range("I50").Value = 39000
$ws.Range("J50").Value = 35632.668
$ws.Range("K50").Value = 39000
$ws.Range("L50").Value = 35632.668
$ws.Range("M50").Value = -38502
$ws.Range("N50").Value = -36628.668

$ws.Range("H102").Value = 11024.643
$ws.Range("I102").Value = 13276
$ws.Range("K102").Value = 13276
$ws.Range("M102").Value = -11654

$ws.Range("H122").Value = 3452.5293
$ws.Range("I122").Value = 2573.76
$ws.Range("K122").Value = 7721.280000000001
$ws.Range("M122").Value = -5271.280000000001

$ws.Range("H132").Value = 2070
$ws.Range("I132").Value = 932.75
$ws.Range("K132").Value = 2798.25
$ws.Range("M132").Value = -268.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 500.81818
$ws.Range("J16").Value = 2994
$ws.Range("K16").Value = 500.81818
$ws.Range("L16").Value = 2994
$ws.Range("M16").Value = -330.81818
$ws.Range("N16").Value = -3334

$ws.Range("H22").Value = 1476.875
$ws.Range("J22").Value = 1504.6666
$ws.Range("L22").Value = 1504.6666
$ws.Range("N22").Value = -2094.6666

$ws.Range("H27").Value = 1476.875
$ws.Range("J27").Value = 1504.6666
$ws.Range("L27").Value = 1504.6666
$ws.Range("N27").Value = -1718.6666

$ws.Range("H31").Value = 13748.5
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 13748.5
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 13748.5
$ws.Range("N31").Value = -14244.5
$ws.Range("M31").ClearContents()

$ws.Range("H40").Value = 2627.1428
$ws.Range("I40").Value = 2198.25
$ws.Range("J40").Value = 3199
$ws.Range("K40").Value = 2198.25
$ws.Range("L40").Value = 3199
$ws.Range("M40").Value = -2062.25
$ws.Range("N40").Value = -3471

$ws.Range("H61").Value = 4188.25
$ws.Range("I61").Value = 4188.25
$ws.Range("K61").Value = 4188.25
$ws.Range("M61").Value = -3986.25

$ws.Range("H113").Value = 4188.25
$ws.Range("I113").Value = 4188.25
$ws.Range("K113").Value = 4188.25
$ws.Range("M113").Value = -2018.25

$ws.Range("H122").Value = 6418.5
$ws.Range("J122").Value = 7108.8
$ws.Range("L122").Value = 21326.4
$ws.Range("N122").Value = -26226.4

$ws.Range("H132").Value = 2678.6072
$ws.Range("I132").Value = 2328.923
$ws.Range("J132").Value = 2981.6667
$ws.Range("K132").Value = 6986.768999999999
$ws.Range("L132").Value = 8945.000100000001
$ws.Range("M132").Value = -4456.768999999999
$ws.Range("N132").Value = -14005.0001

$ws.Range("H136").Value = 2258.524
$ws.Range("I136").Value = 2123.25
$ws.Range("J136").Value = 2438.889
$ws.Range("K136").Value = 6369.75
$ws.Range("L136").Value = 7316.667
$ws.Range("M136").Value = -3819.75
$ws.Range("N136").Value = -12416.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 68017.14
$ws.Range("J46").Value = 68017.14
$ws.Range("L46").Value = 68017.14
$ws.Range("N46").Value = -68479.14

$ws.Range("H126").Value = 2154.3684
$ws.Range("I126").Value = 2275.5454
$ws.Range("J126").Value = 1987.75
$ws.Range("K126").Value = 6826.6362
$ws.Range("L126").Value = 5963.25
$ws.Range("M126").Value = -4356.6362
$ws.Range("N126").Value = -10903.25

$ws.Range("H132").Value = 2924.625
$ws.Range("I132").Value = 3062.8
$ws.Range("J132").Value = 2694.3333
$ws.Range("K132").Value = 9188.400000000001
$ws.Range("L132").Value = 8082.999899999999
$ws.Range("M132").Value = -6658.400000000001
$ws.Range("N132").Value = -13142.9999

$ws.Range("H134").Value = 68017.14
$ws.Range("J134").Value = 68017.14
$ws.Range("L134").Value = 204051.42
$ws.Range("N134").Value = -209121.42
